$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename existing header labels to the lowercase "metadata4ing" convention
$ws.Range("D1").Value = "metadata4ing_IRI"
$ws.Range("E1").Value = "metadata4ing_DESC"

# Add new column F: header (matching the formatting of the other header cells)
# plus the definition value for row 2.
$ws.Range("F1").Value = "metadata4ing_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("F2").Value = "[locstr('A role is the function of an entity or agent with respect to an activity, in the context of a usage, generation, invalidation, association, start, and end.', 'en')]"
